# Add timesheet entries for 11-02-2020 (Feb 11 2020) to the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 101 is a blank separator row (same styling as the other separator rows,
# e.g. row 82 / row 92). Copy its formatting down onto the new blank row.
$ws.Range("A82:C82").Copy()
$ws.Range("A101:C101").PasteSpecial(-4122)

# Rows 102-110 hold the new day's entries. Copy the standard data-row
# formatting (style used throughout the sheet, e.g. row 83) across the new
# range before filling in the values.
$ws.Range("A83:C83").Copy()
$ws.Range("A102:C110").PasteSpecial(-4122)

$entries = @(
  @("Feb 11 10:00 to 11:00", "Uploaded timesheet of previous day to git. Working on making hourly buckets.", "Infimetrics"),
  @("Feb 11 11:00 to 12:00", "Using timedelta adding hour to timestamp data", "Infimetrics"),
  @("Feb 11 12:00 to 13:00", "Created hourly bucket, working on making the algorithm faster for hourly bucket.", "Infimetrics"),
  @("Feb 11 13:00 to 14:00", "Modified logic of creating hourly buckets, veryfied output by writing data in csv file.", "Infimetrics"),
  @("Feb 11 14:00 to 15:00", "Lunch", "Infimetrics"),
  @("Feb 11 15:00 to 16:00", "Logic to create hourly bucket was wrong, working on new logic", "Infimetrics"),
  @("Feb 11 16:00 to 17:00", "Logic to create hourly bucket was wrong, working on new logic", "Infimetrics"),
  @("Feb 11 17:00 to 18:00", "Created hourly bucket, working on connecting other two files.", "Infimetrics"),
  @("Feb 11 18:00 to 19:00", "Working on making generalized code which should work for all three files. ", "Infimetrics")
)

$r = 102
foreach ($entry in $entries) {
  $ws.Cells.Item($r, 1).Value = $entry[0]
  $ws.Cells.Item($r, 2).Value = $entry[1]
  $ws.Cells.Item($r, 3).Value = $entry[2]
  $r = $r + 1
}

# Match the saved view state: scrolled so row 97 is at the top, with D110
# selected as the active cell.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 97
$win.ScrollColumn = 1
$ws.Range("D110").Select()
